$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 (B13/C13 = "5464150 - Mariana Consiglio Kasemodel", no A value)
# was deleted, shifting every row below it up by one (this also brings the row
# heights of rows 14-25 up to rows 13-24 automatically, and shrinks the used
# range from A1:C25 to A1:C24).
$ws.Rows(13).Delete()

# A handful of the B/C (value) cells did not simply follow their row as it
# shifted up -- their text content changed independently. Patch those cells
# to their final values explicitly.
$ws.Range("B10").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C10").Value = "5464150 - Mariana Consiglio Kasemodel"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

$ws.Range("B18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C18").Value = "5464150 - Mariana Consiglio Kasemodel"

$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

$ws.Range("B20").Value = "Média ponderada de provas  e atividades."
$ws.Range("C20").Value = "Média ponderada de provas  e atividades."

$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"
